$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.607.73"
$ws.Range("E2").Value = "  +1.25%  "

$ws.Range("D3").Value = "2.233.43"
$ws.Range("E3").Value = "  -0.07%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "301.85"
$ws.Range("E5").Value = "  +2.41%  "

$ws.Range("D6").Value = "89.57"
$ws.Range("E6").Value = "  +2.19%  "

$ws.Range("E7").Value = "  +0.94%  "

$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  +0.09%  "

$ws.Range("D9").Value = "0.477"
$ws.Range("E9").Value = "  +0.76%  "

$ws.Range("D10").Value = "52.81"
$ws.Range("E10").Value = "  +8.09%  "

$ws.Range("D11").Value = "31.52"
$ws.Range("E11").Value = "  +3.88%  "

$ws.Range("D12").Value = "0.0789"
$ws.Range("E12").Value = "  +0.95%  "

$ws.Range("E13").Value = "  +2.60%  "

$ws.Range("D14").Value = "6.49"
$ws.Range("E14").Value = "  +0.00%  "

$ws.Range("D15").Value = "2.578.27"
$ws.Range("E15").Value = "  +0.13%  "

$ws.Range("D16").Value = "13.98"
$ws.Range("E16").Value = "  +0.99%  "

$ws.Range("D17").Value = "2.227.57"
$ws.Range("E17").Value = "  +0.37%  "

$ws.Range("D18").Value = "0.744"
$ws.Range("E18").Value = "  +1.85%  "

$ws.Range("D19").Value = "40.494.72"
$ws.Range("E19").Value = "  +1.17%  "

$ws.Range("D20").Value = "0.0₃0895"
$ws.Range("E20").Value = "  +0.77%  "

$ws.Range("E21").Value = "  +0.44%  "

$ws.Range("D22").Value = "5.83"
$ws.Range("E22").Value = "  -0.01%  "

$ws.Range("E23").Value = "  +0.54%  "

$ws.Range("D24").Value = "237.52"
$ws.Range("E24").Value = "  +0.39%  "

$ws.Range("D25").Value = "2.54"
$ws.Range("E25").Value = "  +3.14%  "

$ws.Range("E26").Value = "  -0.04%  "

$ws.Range("D27").Value = "1.85"
$ws.Range("E27").Value = "  +1.47%  "

$ws.Range("D28").Value = "23.77"
$ws.Range("E28").Value = "  +4.35%  "

$ws.Range("D30").Value = "9.43"
$ws.Range("E30").Value = "  +2.08%  "

$ws.Range("D31").Value = "157.11"
$ws.Range("E31").Value = "  +0.78%  "

$ws.Range("D32").Value = "32.87"
$ws.Range("E32").Value = "  +2.87%  "

$ws.Range("E33").Value = "  +0.09%  "

$ws.Range("D34").Value = "5.05"
$ws.Range("E34").Value = "  +2.10%  "

$ws.Range("E35").Value = "  +0.68%  "

$ws.Range("D36").Value = "3.01"
$ws.Range("E36").Value = "  +4.39%  "

$ws.Range("D37").Value = "2.34"

$ws.Range("D38").Value = "0.104"
$ws.Range("E38").Value = "  +6.75%  "

$ws.Range("E39").Value = "  +1.71%  "

$ws.Range("D40").Value = "1.76"
$ws.Range("E40").Value = "  +4.44%  "

$ws.Range("D41").Value = "15.84"
$ws.Range("E41").Value = "  +0.61%  "

$ws.Range("D42").Value = "3.86"
$ws.Range("E42").Value = "  -0.03%  "

$ws.Range("D43").Value = "2.087.33"
$ws.Range("E43").Value = "  -1.51%  "

$ws.Range("D44").Value = "19.65"
$ws.Range("E44").Value = "  +9.60%  "

$ws.Range("D45").Value = "0.0274"
$ws.Range("E45").Value = "  +2.44%  "

$ws.Range("D46").Value = "10.11"
$ws.Range("E46").Value = "  +2.53%  "

$ws.Range("D47").Value = "2.86"
$ws.Range("E47").Value = "  +7.11%  "

$ws.Range("D48").Value = "1.86"
$ws.Range("E48").Value = "  -12.71%  "

$ws.Range("D49").Value = "2.441.54"
$ws.Range("E49").Value = "  +0.18%  "

$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").Value = "1.51"
$ws.Range("E50").Value = "  +2.35%  "

$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").Value = "1.14"
$ws.Range("E51").Value = "  +2.59%  "
